# Activities updated for 2022: add a new speaking engagement as the first entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the AMSCUE/Online entry), shifting
# all the existing entries down by one.
$ws.Rows("2:2").Insert()

$ws.Range("A2").Value = 2021
$ws.Range("B2").Value = "Pisa"
$ws.Range("C2").Value = "UNIPI International Workshop on Multidisciplinary studies for sustainable agriculture"
$ws.Range("D2").Value = "Speaker"
$ws.Range("E2").Value = """Evolutionary Biology and Genetic Pest Control"""
$ws.Range("F2").Value = "Research"
